$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table -1.1")
$ws.Activate()

# Q3 answer (C5): Total amount of investment (USD) -> corrected value
$ws.Range("C5").Value = 66368

# Q4 answer (C8): Are there any companies in the rounds2 file which are not
# present in companies? Answer Y/N -> "N"
$ws.Range("C8").Value = "N"

# Leave the selection where the author last left off while editing (C9)
$ws.Range("C9").Select()
